$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "DATA"
$ws.Range("B1").Value = "DAY"
$ws.Range("C1").Value = "TIME SPEND"
$ws.Range("D1").Value = "REMARKS"

$ws.Range("A2").Value = 45159
$ws.Range("A2").NumberFormat = "mm-dd-yy"
$ws.Range("B2").Value = "MONDAY"
$ws.Range("C2").Value = "log-in 10.30 am and log out 6.30"
$ws.Range("D2").Value = "create git and repository(3hrs)"

$ws.Columns.Item(1).ColumnWidth = 9.666666666666666
$ws.Columns.Item(3).ColumnWidth = 27.5
$ws.Columns.Item(4).ColumnWidth = 35.666666666666664

$ws.Range("D2").Select()
